# Daily update of covid19 tracker data files
# The "Country Updates" sheet's Date column (B) was bumped from 43936
# (15/Apr/2020) to 43937 (16/Apr/2020) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

$dateRange = $ws.Range("B5:B96")
foreach ($cell in $dateRange.Cells) {
    if ($cell.Value2 -eq 43936) {
        $cell.Value2 = 43937
    }
}
